$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B26:F26 values from 0 to 2
$ws.Range("B26:F26").Value = 2

# Update the selection in the sheet view to G26
$ws.Range("G26").Select()

$wb.Save()
